$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new data row 21 (new "dailylevelfinish" condition entry) ---
# Copy the formatting of row 13 (same style pattern: s=7 across A:O, s=16 on I)
# into the new row 21 so the cell styles match exactly.
$ws.Range("A13:O13").Copy()
$ws.Range("A21:O21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows(21).RowHeight = 13.5

$ws.Range("A21").Value = "dailylevelfinish"
$ws.Range("B21").Value = "每日关卡完成次数"
$ws.Range("D21").Value = "player"
$ws.Range("E21").Value = "dailylevelfinish"
$ws.Range("F21").Value = 3
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = "3"
$ws.Range("J21").Value = 16
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 1

# --- Conditional formatting: a new "duplicate values" rule on A21, and the
# existing A14 block's range grows to include the newly inserted row. ---
$ws.Range("A14").FormatConditions.Delete()

$fc14 = $ws.Range("A14 A22:A65430").FormatConditions.AddUniqueValues()
$fc14.DupeUnique = 1
$fc14.Font.Color = 10223622
$fc14.Interior.Color = 16762830
$fc14.StopIfTrue = $true

$fc21 = $ws.Range("A21").FormatConditions.AddUniqueValues()
$fc21.DupeUnique = 1
$fc21.Font.Color = 10223622
$fc21.Interior.Color = 16762830
$fc21.StopIfTrue = $true

# --- View state: scrolled/selected cell matches the author's saved view ---
$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J4").Select()
